{"js": "// Apply the billing-letter edits:\n// 1) MEMO-REG-LMD-2020-Aug- -> MEMO-REG-LMD-2020-Sep-\n// 2) August 20, 2020 -> September 14, 2020\n// 3) Cotabato Sugar Central Company, Inc. -> Davao Sugar Central Company, Inc.\n// 4) Brgy. Kilada, Matalam, Cotabato -> 5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\n// 5) Dear Pres. Gotianun: -> Dear Mr. Gotianun:\n// 6) Page margins: left/right 2200 -> 1500 twips (1.5in -> ~1.04in / \"1inch\")\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\"MEMO-REG-LMD-2020-Aug-\", \"MEMO-REG-LMD-2020-Sep-\");\nawait replaceOnce(\"August 20, 2020\", \"September 14, 2020\");\nawait replaceOnce(\"Cotabato Sugar Central Company, Inc.\", \"Davao Sugar Central Company, Inc.\");\nawait replaceOnce(\n  \"Brgy. Kilada, Matalam, Cotabato\",\n  \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\"\n);\nawait replaceOnce(\"Dear Pres. Gotianun:\", \"Dear Mr. Gotianun:\");\n\n// Update left/right page margins from 2200 twips (1.527...in) to 1500 twips.\nconst section = context.document.sections.getFirst();\nconst pageSetup = section.pageSetup;\npageSetup.load(\"leftMargin,rightMargin\");\nawait context.sync();\n\n// 1500 twips / 1440 twips-per-inch = 1.0416666... inch; Office.js page margins are in points (1 twip = 0.05 pt).\nconst newMarginPoints = (1500 / 1440) * 72;\npageSetup.leftMargin = newMarginPoints;\npageSetup.rightMargin = newMarginPoints;\nawait context.sync();\n", "ps1": "# Apply the billing-letter edits:\n# 1) MEMO-REG-LMD-2020-Aug- -> MEMO-REG-LMD-2020-Sep-\n# 2) August 20, 2020 -> September 14, 2020\n# 3) Cotabato Sugar Central Company, Inc. -> Davao Sugar Central Company, Inc.\n# 4) Brgy. Kilada, Matalam, Cotabato -> 5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\n# 5) Dear Pres. Gotianun: -> Dear Mr. Gotianun:\n# 6) Page margins: left/right 2200 -> 1500 twips (1.5in -> ~1.04in / \"1inch\")\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    # Args: FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #       MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n    # MatchWholeWord must stay $false since some search strings end in punctuation\n    # (e.g. a trailing hyphen), which would never be considered a \"whole word\".\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-Text \"MEMO-REG-LMD-2020-Aug-\" \"MEMO-REG-LMD-2020-Sep-\"\nReplace-Text \"August 20, 2020\" \"September 14, 2020\"\nReplace-Text \"Cotabato Sugar Central Company, Inc.\" \"Davao Sugar Central Company, Inc.\"\nReplace-Text \"Brgy. Kilada, Matalam, Cotabato\" \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\"\nReplace-Text \"Dear Pres. Gotianun:\" \"Dear Mr. Gotianun:\"\n\n# Update left/right page margins from 2200 twips (1.527...in) to 1500 twips.\n# COM PageSetup margins are expressed in points (1 inch = 72 points, 1440 twips = 1 inch).\n$newMarginPoints = (1500 / 1440) * 72\n$d.PageSetup.LeftMargin = $newMarginPoints\n$d.PageSetup.RightMargin = $newMarginPoints\n"}
